$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.2452781377806629
$ws.Range("D2").Value = 0.8085146263414069

# Row 3
$ws.Range("C3").Value = 0.01045015380405306
$ws.Range("D3").Value = 0.9917563221715491

# Row 4
$ws.Range("C4").Value = 2.0134224619898
$ws.Range("D4").Value = 0.0564633738587268
$ws.Range("G4").Value = "No"

# Row 5
$ws.Range("C5").Value = 1.884882997438161
$ws.Range("D5").Value = 0.07272435716539905

# Row 6
$ws.Range("C6").Value = -0.2091065006143031
$ws.Range("D6").Value = 0.8362892009256524

# Row 7
$ws.Range("C7").Value = 2.36135410837024
$ws.Range("D7").Value = 0.02748074083366503

# Row 8
$ws.Range("C8").Value = 1.787274704295283
$ws.Range("D8").Value = 0.08767975712171272

# Row 9
$ws.Range("C9").Value = 2.090703067860132
$ws.Range("D9").Value = 0.04832252090076539

# Row 10
$ws.Range("C10").Value = 2.830511613299944
$ws.Range("D10").Value = 0.009735881936968838
$ws.Range("G10").Value = "Sí"

# Row 11
$ws.Range("C11").Value = -0.7314275095761226
$ws.Range("D11").Value = 0.4722376660722789
